# Update leve profit calculations (currentAveragePrice / LevePrice / LeveProfit columns)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, per scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 297
$ws.Range("I2").Value = 131
$ws.Range("J2").Value = 795
$ws.Range("K2").Value = 131
$ws.Range("L2").Value = 795
$ws.Range("M2").Value = -18
$ws.Range("N2").Value = -1021

$ws.Range("H9").Value = 1406.7
$ws.Range("I9").Value = 247
$ws.Range("J9").Value = 2179.8333
$ws.Range("K9").Value = 247
$ws.Range("L9").Value = 2179.8333
$ws.Range("M9").Value = -78
$ws.Range("N9").Value = -2517.8333

$ws.Range("H17").Value = 1431730.2
$ws.Range("J17").Value = 1431730.2
$ws.Range("L17").Value = 4295190.6
$ws.Range("N17").Value = -4295526.6

$ws.Range("H70").Value = 1883.2778
$ws.Range("I70").Value = 1693.2
$ws.Range("J70").Value = 2120.875
$ws.Range("K70").Value = 5079.6
$ws.Range("L70").Value = 6362.625
$ws.Range("M70").Value = -4809.6
$ws.Range("N70").Value = -6902.625

$ws.Range("H73").Value = 1883.2778
$ws.Range("I73").Value = 1693.2
$ws.Range("J73").Value = 2120.875
$ws.Range("K73").Value = 5079.6
$ws.Range("L73").Value = 6362.625
$ws.Range("M73").Value = -4143.6
$ws.Range("N73").Value = -8234.625

$ws.Range("H96").Value = 2447.2
$ws.Range("I96").Value = 2091.6365
$ws.Range("K96").Value = 6274.9095
$ws.Range("M96").Value = -4901.9095

$ws.Range("H100").Value = 1494.75
$ws.Range("I100").Value = 994.5
$ws.Range("K100").Value = 994.5
$ws.Range("M100").Value = -453.5

$ws.Range("H107").Value = 1160.5
$ws.Range("I107").Value = 1051.4286
$ws.Range("J107").Value = 1415
$ws.Range("K107").Value = 1051.4286
$ws.Range("L107").Value = 1415
$ws.Range("M107").Value = 868.5714
$ws.Range("N107").Value = -5255

$ws.Range("H116").Value = 3989.6667
$ws.Range("I116").Value = 3487.3333
$ws.Range("J116").Value = 5496.6665
$ws.Range("K116").Value = 3487.3333
$ws.Range("L116").Value = 5496.6665
$ws.Range("M116").Value = -45.33329999999978
$ws.Range("N116").Value = -12380.6665

$ws.Range("H129").Value = 18913.076
$ws.Range("I129").Value = 32348.773
$ws.Range("K129").Value = 97046.319
$ws.Range("M129").Value = -92046.319

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 300669.4
$ws.Range("I13").Value = 300669.4
$ws.Range("K13").Value = 300669.4
$ws.Range("M13").Value = -300525.4

$ws.Range("H61").Value = 634727.4399999999
$ws.Range("I61").Value = 2641.3333
$ws.Range("J61").Value = 2703372.8
$ws.Range("K61").Value = 2641.3333
$ws.Range("L61").Value = 2703372.8
$ws.Range("M61").Value = -2429.3333
$ws.Range("N61").Value = -2703796.8

$ws.Range("H74").Value = 14492.105
$ws.Range("I74").Value = 1038.0834
$ws.Range("J74").Value = 37556.145
$ws.Range("K74").Value = 1038.0834
$ws.Range("L74").Value = 37556.145
$ws.Range("M74").Value = -164.0834
$ws.Range("N74").Value = -39304.145

$ws.Range("H77").Value = 14492.105
$ws.Range("I77").Value = 1038.0834
$ws.Range("J77").Value = 37556.145
$ws.Range("K77").Value = 5190.416999999999
$ws.Range("L77").Value = 187780.725
$ws.Range("M77").Value = -822.4169999999995
$ws.Range("N77").Value = -196516.725

$ws.Range("H110").Value = 5480.3076
$ws.Range("I110").Value = 6383.8096
$ws.Range("K110").Value = 6383.8096
$ws.Range("M110").Value = -4338.8096

$ws.Range("H136").Value = 634727.4399999999
$ws.Range("I136").Value = 2641.3333
$ws.Range("J136").Value = 2703372.8
$ws.Range("K136").Value = 7923.999899999999
$ws.Range("L136").Value = 8110118.399999999
$ws.Range("M136").Value = -5373.999899999999
$ws.Range("N136").Value = -8115218.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1194.9375
$ws.Range("I105").Value = 1415.3334
$ws.Range("K105").Value = 1415.3334
$ws.Range("M105").Value = 331.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 6000
$ws.Range("I45").Value = 6000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 6000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -5407
$ws.Range("N45").ClearContents()

$ws.Range("H58").Value = 19895.904
$ws.Range("I58").Value = 7442.6665
$ws.Range("K58").Value = 7442.6665
$ws.Range("M58").Value = -7239.6665

$ws.Range("H62").Value = 6345.5454
$ws.Range("I62").Value = 6186
$ws.Range("J62").Value = 6624.75
$ws.Range("K62").Value = 6186
$ws.Range("L62").Value = 6624.75
$ws.Range("M62").Value = -5562
$ws.Range("N62").Value = -7872.75

$ws.Range("H65").Value = 6345.5454
$ws.Range("I65").Value = 6186
$ws.Range("J65").Value = 6624.75
$ws.Range("K65").Value = 30930
$ws.Range("L65").Value = 33123.75
$ws.Range("M65").Value = -27810
$ws.Range("N65").Value = -39363.75

$ws.Range("H86").Value = 10547.789
$ws.Range("I86").Value = 11053.588
$ws.Range("K86").Value = 11053.588
$ws.Range("M86").Value = -9930.588

$ws.Range("H89").Value = 10547.789
$ws.Range("I89").Value = 11053.588
$ws.Range("K89").Value = 55267.94
$ws.Range("M89").Value = -49651.94

$ws.Range("H99").Value = 5280
$ws.Range("J99").Value = 3750
$ws.Range("L99").Value = 3750
$ws.Range("N99").Value = -6746

$ws.Range("H122").Value = 2407.3076
$ws.Range("I122").Value = 1412.75
$ws.Range("K122").Value = 4238.25
$ws.Range("M122").Value = -1788.25

$ws.Range("H126").Value = 5280
$ws.Range("J126").Value = 3750
$ws.Range("L126").Value = 11250
$ws.Range("N126").Value = -16190

$ws.Range("H136").Value = 19895.904
$ws.Range("I136").Value = 7442.6665
$ws.Range("K136").Value = 22327.9995
$ws.Range("M136").Value = -19777.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 47.75
$ws.Range("I17").Value = 28.4
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 85.19999999999999
$ws.Range("L17").Value = 240
$ws.Range("M17").Value = 83.80000000000001
$ws.Range("N17").Value = -578

$ws.Range("H68").Value = 1753.2084
$ws.Range("J68").Value = 1396.5333
$ws.Range("L68").Value = 4189.5999
$ws.Range("N68").Value = -5811.5999

$ws.Range("H71").Value = 1753.2084
$ws.Range("J71").Value = 1396.5333
$ws.Range("L71").Value = 12568.7997
$ws.Range("N71").Value = -20680.7997

$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -754
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 2183.1667
$ws.Range("J107").Value = 2449.75
$ws.Range("L107").Value = 7349.25
$ws.Range("N107").Value = -11189.25

$ws.Range("H113").Value = 1039.6
$ws.Range("J113").Value = 1059.4546
$ws.Range("L113").Value = 3178.3638
$ws.Range("N113").Value = -7518.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6020.2
$ws.Range("J11").Value = 15080.8
$ws.Range("L11").Value = 15080.8
$ws.Range("N11").Value = -15358.8

$ws.Range("H97").Value = 1050.6389
$ws.Range("I97").Value = 857.0417
$ws.Range("K97").Value = 857.0417
$ws.Range("M97").Value = -361.0417

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 30000
$ws.Range("I7").Value = 50000
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 50000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -49888
$ws.Range("N7").Value = -10224

$ws.Range("H22").Value = 1554.8422
$ws.Range("J22").Value = 2584.5715
$ws.Range("L22").Value = 2584.5715
$ws.Range("N22").Value = -3174.5715

$ws.Range("H27").Value = 1554.8422
$ws.Range("J27").Value = 2584.5715
$ws.Range("L27").Value = 2584.5715
$ws.Range("N27").Value = -2798.5715

$ws.Range("H55").Value = 2239.111
$ws.Range("I55").Value = 2356.2856
$ws.Range("J55").Value = 2164.5454
$ws.Range("K55").Value = 2356.2856
$ws.Range("L55").Value = 2164.5454
$ws.Range("M55").Value = -2183.2856
$ws.Range("N55").Value = -2510.5454

$ws.Range("H93").Value = 7329.5
$ws.Range("I93").Value = 25675
$ws.Range("J93").Value = 1214.3334
$ws.Range("K93").Value = 25675
$ws.Range("L93").Value = 1214.3334
$ws.Range("M93").Value = -24427
$ws.Range("N93").Value = -3710.3334

$ws.Range("H100").Value = 2792.0454
$ws.Range("I100").Value = 2518.3076
$ws.Range("J100").Value = 3187.4443
$ws.Range("K100").Value = 2518.3076
$ws.Range("L100").Value = 3187.4443
$ws.Range("M100").Value = -1977.3076
$ws.Range("N100").Value = -4269.4443

$ws.Range("H122").Value = 5404.5
$ws.Range("I122").Value = 4824.5
$ws.Range("K122").Value = 14473.5
$ws.Range("M122").Value = -12023.5

$ws.Range("H126").Value = 30000
$ws.Range("I126").Value = 50000
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 150000
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -147530
$ws.Range("N126").Value = -34940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 84474.8
$ws.Range("I9").Value = 93093.5
$ws.Range("K9").Value = 93093.5
$ws.Range("M9").Value = -92953.5

$ws.Range("H31").Value = 13708.333
$ws.Range("I31").Value = 9250
$ws.Range("K31").Value = 9250
$ws.Range("M31").Value = -8902

$ws.Range("H47").Value = 52000
$ws.Range("I47").Value = 52000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 52000
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -51428
$ws.Range("N47").ClearContents()

$ws.Range("H126").Value = 42193.11
$ws.Range("I126").Value = 42193.11
$ws.Range("K126").Value = 126579.33
$ws.Range("M126").Value = -124109.33
